$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 214, shifting existing rows 214-222 down to 215-223.
$ws.Rows(214).Insert()

# Populate the newly inserted row 214 with the new data record.
$ws.Range("A214").Value = 6
$ws.Range("B214").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C214").Value = "Metropolitana"
$ws.Range("D214").Value = 44753
$ws.Range("E214").Value = 13
$ws.Range("F214").Value = 100112001
$ws.Range("G214").Value = "Berenjena"
$ws.Range("H214").Value = "Sin especificar"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 400
$ws.Range("K214").Value = 9000
$ws.Range("L214").Value = 10000
$ws.Range("M214").Value = 9425
$ws.Range("N214").Value = "$/caja 50 unidades"
$ws.Range("O214").Value = "Región de Arica y Parinacota"
$ws.Range("P214").Value = 188
$ws.Range("Q214").Value = 50
$ws.Range("R214").Value = "Hortaliza"
